# Lattice-multiplication exercise refresh: replace the contents of every
# cell in the single 5x3 practice table with the new set of problems.
#
# Each cell's text is five "lines" (separated by manual line breaks, i.e.
# a vertical-tab char in the Range.Text model):
#   1) "AA x BB"          - the two factors
#   2) "  C    D"         - the digits of the second factor, spaced out
#   3) "  ----"           - separator
#   4) "E|    |"          - first digit of first factor
#   5) "F|    |"          - second digit of first factor

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vbr = [char]11

$newCells = @(
    @("81 x 96", "  9    6", "8|    |", "1|    |"),
    @("25 x 89", "  8    9", "2|    |", "5|    |"),
    @("22 x 27", "  2    7", "2|    |", "2|    |"),
    @("30 x 88", "  8    8", "3|    |", "0|    |"),
    @("27 x 63", "  6    3", "2|    |", "7|    |"),
    @("74 x 69", "  6    9", "7|    |", "4|    |"),
    @("73 x 87", "  8    7", "7|    |", "3|    |"),
    @("61 x 61", "  6    1", "6|    |", "1|    |"),
    @("61 x 22", "  2    2", "6|    |", "1|    |"),
    @("23 x 10", "  1    0", "2|    |", "3|    |"),
    @("92 x 28", "  2    8", "9|    |", "2|    |"),
    @("84 x 36", "  3    6", "8|    |", "4|    |"),
    @("65 x 26", "  2    6", "6|    |", "5|    |"),
    @("34 x 83", "  8    3", "3|    |", "4|    |"),
    @("11 x 88", "  8    8", "1|    |", "1|    |")
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $vals = $newCells[$idx]
        $idx = $idx + 1

        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $vals[0] + $vbr + $vals[1] + $vbr + "  ----" + $vbr + $vals[2] + $vbr + $vals[3]
    }
}

Write-Output "done: updated $idx cells"
